# Applies the commit "Added header with case management location":
# In the document header paragraph that reads
#   <<cs_{writtenByJudge}>><<hearingLocation.venue_name>><<else>> Online Civil Claims<<es_>>
# the placeholder "hearingLocation.venue_name" becomes "caseManagementLocation.venue_name".
#
# We use a narrow, unambiguous Find/Replace on the unique substring
# "hearingLocation.venue_name" so the other (unrelated) occurrences of
# "hearingLocation" elsewhere in the template (e.g. hearingLocation.site_name,
# hearingLocation.court_address, hearingLocation.postcode) are left untouched.

$d = $word.ActiveDocument

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute(
    "hearingLocation.venue_name",  # FindText
    $true,                          # MatchCase
    $true,                          # MatchWholeWord
    $false,                         # MatchWildcards
    $false,                         # MatchSoundsLike
    $false,                         # MatchAllWordForms
    $true,                          # Forward
    1,                               # Wrap (wdFindContinue)
    $false,                         # Format
    "caseManagementLocation.venue_name", # ReplaceWith
    2                                # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find target text 'hearingLocation.venue_name' to replace."
}

$d.Save()
